$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing row (row 6) down into the new row 7
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)

# Populate the new row with the announcement test case data
$ws.Range("A7").Value = "EXL_CorporateLensHomePage_Announcements"
$ws.Range("B7").Value = "Add Announcement"
$ws.Range("C7").Value = "N"
$ws.Range("D7").Value = "Y"
$ws.Range("F7").Value = "Sprint1"

# Extend the data validations to cover the newly added row
$ws.Range("C2:D6").Validation.Delete()
$ws.Range("C2:D7").Validation.Add(3, 1, 1, '"Y,N"')

$ws.Range("F2:F6").Validation.Delete()
$ws.Range("F2:F7").Validation.Add(3, 1, 1, '"Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10"')

# Update the selection/active cell to the newly added row, matching the marked-for-execution state
$ws.Range("A7").Select()
